$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Replace the combined (3-item) summary paragraph FIRST.
# Its full text is a unique string in the document, but the three
# individual numbered sentences it is built from also appear (with
# different numbering) as separate <w:t> runs earlier in the document.
# Handling this exact, full, unique string before touching the shorter
# per-item strings avoids any accidental double-replacement.
# ------------------------------------------------------------------
$old3 = "1.球阀表面水平开设介质穿过的流通孔，阀芯通过顶部阀杆转动控制流通孔转向。 2.连筒和密封环数量均为两个，护筒两侧焊接固定阀体的法兰盘，连筒位于法兰盘内部中心。 3.护套套接在护筒表面与阀杆的连接处，护套防止介质从通孔流出，阀杆在护套内沿通孔水平转动，阀芯通过转轴在卡板底面转动。"
$new3 = "1.护套套接在护筒表面与阀杆的连接处，护套防止介质从通孔流出，阀杆在护套内沿通孔水平转动，阀芯通过转轴在卡板底面转动。 2.阀杆底部通过螺钉与阀芯相连，固定盘顶部安装在护套底面，固定盘底面通过螺钉安装护筒顶部，护筒顶部对应阀杆开设通孔，固定盘底面通过螺钉安装在通孔四周，护筒顶部内壁对应卡板开设方形插槽，阀芯通过卡板水平插接插槽安装在护筒内的通孔底部。 3.连筒和密封环数量均为两个，护筒两侧焊接固定阀体的法兰盘，连筒位于法兰盘内部中心。"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ------------------------------------------------------------------
# Step 2: Rotate the four numbered "核心技术" items (each its own
# <w:t> inside a shared run, separated by <w:br/>). The content that
# used to sit in position 3 now belongs in position 1, the content
# that used to sit in position 4 now belongs in position 2, the
# content that used to sit in position 2 now belongs in position 3,
# and the content that used to sit in position 1 now belongs in
# position 4 (all renumbered 1..4 in their new slot).
#
# A plain Find/Replace only rewrites text *in place*, so to actually
# move content between the four slots we first tag each slot with a
# unique, never-seen-before placeholder token, and only then rewrite
# each placeholder with the final (renumbered) text that belongs
# there. That two-pass approach avoids any slot being overwritten
# before its original content has been captured elsewhere.
# ------------------------------------------------------------------
$slot1Old = "1.球阀表面水平开设介质穿过的流通孔，阀芯通过顶部阀杆转动控制流通孔转向。"
$slot2Old = "2.连筒和密封环数量均为两个，护筒两侧焊接固定阀体的法兰盘，连筒位于法兰盘内部中心。"
$slot3Old = "3.护套套接在护筒表面与阀杆的连接处，护套防止介质从通孔流出，阀杆在护套内沿通孔水平转动，阀芯通过转轴在卡板底面转动。"
$slot4Old = "4.阀杆底部通过螺钉与阀芯相连，固定盘顶部安装在护套底面，固定盘底面通过螺钉安装护筒顶部，护筒顶部对应阀杆开设通孔，固定盘底面通过螺钉安装在通孔四周，护筒顶部内壁对应卡板开设方形插槽，阀芯通过卡板水平插接插槽安装在护筒内的通孔底部。"

$d.Content.Find.Execute($slot1Old, $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT1@@", 2) | Out-Null
$d.Content.Find.Execute($slot2Old, $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT2@@", 2) | Out-Null
$d.Content.Find.Execute($slot3Old, $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT3@@", 2) | Out-Null
$d.Content.Find.Execute($slot4Old, $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT4@@", 2) | Out-Null

$slot1New = "1.护套套接在护筒表面与阀杆的连接处，护套防止介质从通孔流出，阀杆在护套内沿通孔水平转动，阀芯通过转轴在卡板底面转动。"
$slot2New = "2.阀杆底部通过螺钉与阀芯相连，固定盘顶部安装在护套底面，固定盘底面通过螺钉安装护筒顶部，护筒顶部对应阀杆开设通孔，固定盘底面通过螺钉安装在通孔四周，护筒顶部内壁对应卡板开设方形插槽，阀芯通过卡板水平插接插槽安装在护筒内的通孔底部。"
$slot3New = "3.连筒和密封环数量均为两个，护筒两侧焊接固定阀体的法兰盘，连筒位于法兰盘内部中心。"
$slot4New = "4.球阀表面水平开设介质穿过的流通孔，阀芯通过顶部阀杆转动控制流通孔转向。"

$d.Content.Find.Execute("@@SLOT1@@", $true, $false, $false, $false, $false, $true, 1, $false, $slot1New, 2) | Out-Null
$d.Content.Find.Execute("@@SLOT2@@", $true, $false, $false, $false, $false, $true, 1, $false, $slot2New, 2) | Out-Null
$d.Content.Find.Execute("@@SLOT3@@", $true, $false, $false, $false, $false, $true, 1, $false, $slot3New, 2) | Out-Null
$d.Content.Find.Execute("@@SLOT4@@", $true, $false, $false, $false, $false, $true, 1, $false, $slot4New, 2) | Out-Null

# ------------------------------------------------------------------
# Step 3: Extend the two "创新点" sentences with extra trailing text /
# fix final punctuation.
# ------------------------------------------------------------------
$old2a = "1.本便于拆卸的细长密封浮动球阀结构巧妙合理，"
$new2a = "1.本便于拆卸的细长密封浮动球阀结构巧妙合理，控制效果好，便于工人操作。"
$d.Content.Find.Execute($old2a, $true, $false, $false, $false, $false, $true, 1, $false, $new2a, 2) | Out-Null

$old2b = "2.使用细长密封浮动球阀时，通过护筒顶部的方形插槽，将阀芯通过卡板水平插接插槽，安装在护筒内的通孔底部，阀杆底面与阀芯顶部连接，阀芯通过顶部阀杆转动控制流通孔转向，"
$new2b = "2.使用细长密封浮动球阀时，通过护筒顶部的方形插槽，将阀芯通过卡板水平插接插槽，安装在护筒内的通孔底部，阀杆底面与阀芯顶部连接，阀芯通过顶部阀杆转动控制流通孔转向。"
$d.Content.Find.Execute($old2b, $true, $false, $false, $false, $false, $true, 1, $false, $new2b, 2) | Out-Null
